# Updated cryptos list on Thu Aug 29 23:36:02 UTC 2024 with GitHub Actions
# Applies Price (D) and Volume(1h) (E) column updates for rows 2-51 on Sheet1.
# Cells are forced to text (NumberFormat "@") before assignment so that
# purely-numeric-looking strings (e.g. "536.91") are not auto-converted to
# numbers by Excel; the style is reset back to "Normal" afterward so no
# stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 2; D = "59.232.55"; E = "  +0.58%  " },
    @{ Row = 3; D = "2.520.75"; E = "  -0.22%  " },
    @{ Row = 4; D = $null; E = "  +0.02%  " },
    @{ Row = 5; D = "536.91"; E = "  -0.04%  " },
    @{ Row = 6; D = "139.66"; E = "  -2.99%  " },
    @{ Row = 7; D = $null; E = "  +0.01%  " },
    @{ Row = 8; D = $null; E = "  -1.63%  " },
    @{ Row = 9; D = "2.527.68"; E = "  +0.11%  " },
    @{ Row = 10; D = $null; E = "  +0.63%  " },
    @{ Row = 11; D = $null; E = "  +1.25%  " },
    @{ Row = 12; D = "5.45"; E = "  -1.21%  " },
    @{ Row = 13; D = $null; E = "  +1.28%  " },
    @{ Row = 14; D = "2.965.89"; E = "  -0.05%  " },
    @{ Row = 15; D = "59.161.62"; E = "  +0.46%  " },
    @{ Row = 16; D = "22.85"; E = "  -2.83%  " },
    @{ Row = 17; D = $null; E = "  +1.34%  " },
    @{ Row = 18; D = "2.543.96"; E = "  +1.03%  " },
    @{ Row = 19; D = $null; E = "  -2.35%  " },
    @{ Row = 20; D = "4.23"; E = "  -0.60%  " },
    @{ Row = 21; D = "321.81"; E = "  -0.29%  " },
    @{ Row = 22; D = "0.999"; E = "  -0.01%  " },
    @{ Row = 23; D = $null; E = "  +1.18%  " },
    @{ Row = 24; D = "62.14"; E = "  +0.66%  " },
    @{ Row = 25; D = "0.424"; E = "  -2.72%  " },
    @{ Row = 26; D = $null; E = "  +1.48%  " },
    @{ Row = 27; D = "0.999"; E = "  +0.36%  " },
    @{ Row = 28; D = $null; E = "  +0.22%  " },
    @{ Row = 29; D = "6.75"; E = "  +0.53%  " },
    @{ Row = 30; D = $null; E = "  +0.06%  " },
    @{ Row = 31; D = "0.0₃0766"; E = "  -0.27%  " },
    @{ Row = 32; D = "160.76"; E = "  +1.64%  " },
    @{ Row = 33; D = $null; E = "  +0.26%  " },
    @{ Row = 34; D = $null; E = "  +2.04%  " },
    @{ Row = 35; D = $null; E = "  -4.80%  " },
    @{ Row = 36; D = "18.49"; E = "  -0.46%  " },
    @{ Row = 37; D = "4.21"; E = "  -3.08%  " },
    @{ Row = 38; D = $null; E = "  -2.01%  " },
    @{ Row = 39; D = "36.98"; E = "  +1.46%  " },
    @{ Row = 40; D = "3.63"; E = "  -0.14%  " },
    @{ Row = 41; D = "0.804"; E = $null },
    @{ Row = 42; D = "283.69"; E = "  -4.12%  " },
    @{ Row = 43; D = $null; E = "  -6.65%  " },
    @{ Row = 44; D = $null; E = "  +0.06%  " },
    @{ Row = 45; D = $null; E = "  +0.78%  " },
    @{ Row = 46; D = $null; E = "  -1.30%  " },
    @{ Row = 47; D = "0.0930"; E = $null },
    @{ Row = 48; D = "122.48"; E = "  -1.48%  " },
    @{ Row = 49; D = "18.53"; E = "  -0.39%  " },
    @{ Row = 50; D = $null; E = "  -0.36%  " },
    @{ Row = 51; D = $null; E = "  -1.74%  " }
)

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        Set-TextValue ($ws.Cells.Item($u.Row, 4)) $u.D
    }
    if ($null -ne $u.E) {
        Set-TextValue ($ws.Cells.Item($u.Row, 5)) $u.E
    }
}
